# clase_8.pptx - "se actualizo la informacion"
#
# Slide 6 (the "PRACTICA" exercise slide) had its biometric-clock report
# lines re-typed: each of the four "5/7/2015 h:mm:ss AM/PM" lines gained a
# leading space before the date, and a couple of other lines got touched
# during a proof-read pass (the misspelled word "salio" and the second
# date in exercise 4 were re-entered). None of the wording itself changed;
# only exactly where PowerPoint split the runs changed as a side effect of
# the in-place retyping.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Paragraph 5..8 hold the four clock-report lines, e.g. "5/7/2015 9:01:13 AM".
# Re-typing the date added a leading space in front of it, splitting the
# run into " 5/7/2015 " + "9:01:13 AM".
$dateLines = 5, 6, 7, 8
foreach ($idx in $dateLines) {
    $para = $tr.Paragraphs($idx, 1)
    $datePart = $para.Characters(1, 9)
    $datePart.Text = " 5/7/2015 "
}

# Paragraph 12 is exercise 3): the word "salio" was retyped (missing the
# accent on "salió"), which split the sentence into three runs.
$p12 = $tr.Paragraphs(12, 1)
$salioStart = $p12.Text.IndexOf("salio") + 1
$salioRun = $p12.Characters($salioStart, 5)
$salioRun.Text = "salio"

# Paragraph 14 is exercise 4): the trailing date "2001-04-04?" was retyped,
# splitting it into its own run.
$p14 = $tr.Paragraphs(14, 1)
$secondDateStart = $p14.Text.IndexOf("2001-04-04?") + 1
$secondDateRun = $p14.Characters($secondDateStart, 11)
$secondDateRun.Text = "2001-04-04?"
